$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G7").Value = "Sailen"
$ws.Range("G8").Value = "Das"

$ws.Range("G8").Select()
